{"js": "// Update the heading, drop the subtitle paragraph, and turn the materials\n// list into a bulleted list with shortened labels.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Map of the exact original paragraph text -> new paragraph text.\nconst replacements = new Map([\n  [\n    \"Implante de Marcapasso Fisiol\u00f3gico (His/CSP \u2013 Astra\u2122)\",\n    \"Marcapasso Fisiol\u00f3gico (Astra His/CSP)\"\n  ],\n  [\"Gerador \u2013 Astra\u2122\", \"\u2022 Gerador Astra\"],\n  [\"Bainha His \u2013 C315\u2122\", \"\u2022 Bainha C315\"],\n  [\"Eletrodo His / Septal \u2013 3830\", \"\u2022 Eletrodo 3830\"],\n  [\"Eletrodo Atrial \u2013 5076-52\", \"\u2022 Eletrodo 5076-52\"],\n  [\"Ferramenta de Corte\", \"\u2022 Ferramenta de corte\"],\n  [\"Fio Guia\", \"\u2022 Fio guia\"],\n  [\"Introdutor \u2013 2\", \"\u2022 Introdutor \u2013 2\"]\n]);\n\n// The standalone subtitle paragraph is removed entirely.\nconst textToDelete = \"Estimula\u00e7\u00e3o de feixe de His / septal.\";\n\n// Walk from the bottom up so deleting a paragraph doesn't disturb the\n// indices of paragraphs we haven't processed yet.\nfor (let i = items.length - 1; i >= 0; i--) {\n  const para = items[i];\n  const text = para.text;\n  if (text === textToDelete) {\n    para.delete();\n    continue;\n  }\n  const newText = replacements.get(text);\n  if (newText !== undefined) {\n    para.insertText(newText, \"Replace\");\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the heading, drop the subtitle paragraph, and turn the materials\n# list into a bulleted list with shortened labels.\n$d = $word.ActiveDocument\n\n# Map of the exact original paragraph text -> new paragraph text.\n$replacements = @{\n    \"Implante de Marcapasso Fisiol\u00f3gico (His/CSP \u2013 Astra\u2122)\" = \"Marcapasso Fisiol\u00f3gico (Astra His/CSP)\"\n    \"Gerador \u2013 Astra\u2122\"              = \"\u2022 Gerador Astra\"\n    \"Bainha His \u2013 C315\u2122\"            = \"\u2022 Bainha C315\"\n    \"Eletrodo His / Septal \u2013 3830\"  = \"\u2022 Eletrodo 3830\"\n    \"Eletrodo Atrial \u2013 5076-52\"     = \"\u2022 Eletrodo 5076-52\"\n    \"Ferramenta de Corte\"           = \"\u2022 Ferramenta de corte\"\n    \"Fio Guia\"                      = \"\u2022 Fio guia\"\n    \"Introdutor \u2013 2\"                = \"\u2022 Introdutor \u2013 2\"\n}\n\n# The standalone subtitle paragraph is removed entirely.\n$textToDelete = \"Estimula\u00e7\u00e3o de feixe de His / septal.\"\n\n# Walk from the bottom up so deleting a paragraph doesn't disturb the\n# indices of paragraphs we haven't processed yet.\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs.Item($i)\n    $text = $p.Range.Text.TrimEnd(\"`r\", \"`a\")\n    if ($text -eq $textToDelete) {\n        $p.Range.Delete()\n        continue\n    }\n    if ($replacements.ContainsKey($text)) {\n        $p.Range.Text = $replacements[$text]\n    }\n}\n"}
